# Update the "Fitness" values in column C (rows 2-44) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    11370, # C2
    11370, # C3
    11370, # C4
    10511, # C5
    10511, # C6
    10511, # C7
    10365, # C8
    10365, # C9
    10056, # C10
    10056, # C11
    10056, # C12
    10056, # C13
    10056, # C14
    8370,  # C15
    8370,  # C16
    8370,  # C17
    8370,  # C18
    8370,  # C19
    8370,  # C20
    8370,  # C21
    8370,  # C22
    8370,  # C23
    8370,  # C24
    8370,  # C25
    8370,  # C26
    8370,  # C27
    8370,  # C28
    8370,  # C29
    8370,  # C30
    7945,  # C31
    7657,  # C32
    7657,  # C33
    7657,  # C34
    7657,  # C35
    7657,  # C36
    7657,  # C37
    7657,  # C38
    7657,  # C39
    7657,  # C40
    7657,  # C41
    7657,  # C42
    7293,  # C43
    7293   # C44
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
